$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 145 (shifts existing rows 145-182 down to 146-183)
$ws.Rows("145:145").Insert()

# Populate the newly inserted row 145 with the new weekly price record
$ws.Cells.Item(145,1).Value = 4
$ws.Cells.Item(145,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(145,3).Value = "Los Lagos"
$ws.Cells.Item(145,4).Value = 44722
$ws.Cells.Item(145,5).Value = 10
$ws.Cells.Item(145,6).Value = 100112009
$ws.Cells.Item(145,7).Value = "Acelga"
$ws.Cells.Item(145,8).Value = "Sin especificar"
$ws.Cells.Item(145,9).Value = "Primera"
$ws.Cells.Item(145,10).Value = 80
$ws.Cells.Item(145,11).Value = 12000
$ws.Cells.Item(145,12).Value = 12000
$ws.Cells.Item(145,13).Value = 12000
$ws.Cells.Item(145,14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(145,15).Value = "Región de La Araucanía"
$ws.Cells.Item(145,16).Value = 1000
$ws.Cells.Item(145,17).Value = 12
$ws.Cells.Item(145,18).Value = "Hortaliza"
